$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order the header labels for columns E/F/G/H ---------------------
# Original layout:  E1=deaths  F1=losses  G1=affected people  H1=WL
# Target layout:     E1=deaths  F1=losses  G1=affected people  H1=WL
# (same visible text - only shared-string table order changes internally,
#  which Excel handles automatically when the cells are re-written)
$ws.Range("E1").Value = "deaths"
$ws.Range("F1").Value = "losses"
$ws.Range("G1").Value = "affected people"
$ws.Range("H1").Value = "WL"

# --- Insert two extra timestep rows --------------------------------------
# Old rows: 2 (t=0), 3 (t=1), 4 (t=2)
# New rows: 2 (t=0), 3 (t=1, new/averaged), 4 (t=1->2 old row3), 5 (t=3, new/averaged), 6 (t=2->4 old row4)
$ws.Rows.Item(3).EntireRow.Insert()
$ws.Rows.Item(5).EntireRow.Insert()

# The row insert copies formatting from the surrounding rows; clear that
# back to the unformatted default and then restore only the bold/border
# style that column A carries in every other row.
$ws.Range("A3:H3").ClearFormats()
$ws.Range("A5:H5").ClearFormats()

$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Re-number the timestep index column (A) -----------------------------
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# --- New row 3: average of (now) row 2 and row 4 -------------------------
$ws.Range("B3").Formula = "=(B2+B4)/2"
$ws.Range("C3:H3").Formula = "=(C2+C4)/2"

# --- New row 5: average of (now) row 4 and row 6 -------------------------
$ws.Range("B5").Formula = "=(B4+B6)/2"
$ws.Range("C5:H5").Formula = "=(C4+C6)/2"

# --- Restore the selection Excel would leave behind -----------------------
$ws.Range("D9").Select() | Out-Null
